# Append new scrape results to the "ランサーズ" sheet and refresh the
# "captured at" timestamp on every still-present row.
#
# Before: rows 2-5 hold 4 listings, all timestamped 2025-10-13 06:37:36.
# After : rows 2-8 hold 7 listings, all timestamped 2025-10-13 12:38:17 -
#         3 brand-new listings were merged in (sorted by score, column G,
#         descending), two existing ones shifted down one slot, and the
#         last existing one shifted down two slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-10-13 12:38:17"

# Final row data, in sheet order (row 2 .. row 8). $null marks a column
# that must stay completely empty (no cell at all), matching rows that
# never received a skill-summary tag from the scraper.
$rows = @(
    @{ B = "【GAS開発】配送状況管理の自動化を依頼します"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412306"; G = 170; H = "◆開発,自動化 ◇管理" },
    @{ B = "【急募】クリニック向け内視鏡画像システム開発の依頼"; D = "300,000 円 ~ 500,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412233"; G = 125; H = "◆開発,システム開発" },
    @{ B = "【急募】onedrive上のexcelで自動化システム構築依頼"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412194"; G = 95; H = "◆自動化" },
    @{ B = "【急募】スタートアップ向けプロダクト開発のパートナー募集"; D = "300,000 円 ~ 500,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412179"; G = 75; H = "◆開発" },
    @{ B = "【音声コマンド起動】超小型・低電力レコーダーのプロトタイプ開発"; D = "10,000 円 ~ 20,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412261"; G = 60; H = "◆開発" },
    @{ B = "微生物の特定と分類を行いたく、画像解析の専門家を探しています!(急いでません!)"; D = "50,000 円 ~ 100,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5411887"; G = 18; H = $null },
    @{ B = "LINE公式(Lステップ)のリッチメニューの構築"; D = "5,000 円 ~ 10,000 円 / 固定"; F = "https://www.lancers.jp/work/detail/5412357"; G = 10; H = $null }
)

# All cells in every data row share the same category / deadline text.
$category = "システム開発"
$deadline = "期限情報なし"

# Drop every existing hyperlink (and its relationship) up front - row
# positions are about to change and we rebuild F2:F8 from scratch below,
# so there is nothing worth preserving in place.
$ws.Cells.Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $category
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $deadline
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    }

    # Re-create the hyperlink over the URL cell right away, reusing the
    # literal URL we just wrote (reading .Value back from the cell is
    # unreliable in this host, so we never round-trip through it).
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F) | Out-Null
    $ws.Cells.Item($r, 6).Style = "Hyperlink"

    $r = $r + 1
}

# Column H widened from 12 to 13 characters to fit the new skill tags.
$ws.Columns(8).ColumnWidth = 12.15
